$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column F entirely (humanitarian_scope_pct), shifting
# humanitarian_scope_desc_eng (was G) to F and
# humanitarian_scope_desc_fr (was H) to G.
$ws.Range("F1").EntireColumn.Delete()
